$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (column D) and Volume(1h) (column E) figures for the
# crypto list. Several Price values look like plain numbers (e.g. "225.14",
# "1.00", "0.0610") but must stay as literal text, matching how the sheet
# originally stored them (inline/shared strings, not numeric cells). To stop
# Excel from auto-converting these to numbers we temporarily force the cell
# to Text format before assigning the value, then restore the default
# "Normal" style afterwards so no stray formatting is left behind.

$ws.Range("D2").Value = "34.763.99"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.810.61"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.603"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.36"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.68%  "
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0671"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0999"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.41%  "
$ws.Range("D12").Value = "2.074.25"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "1.813.48"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.635"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "34.798.20"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.03%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.121"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0515"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.650"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0188"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +5.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "82.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("E46").Value = "  +3.97%  "
$ws.Range("D47").Value = "1.975.65"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.24%  "

# Rows 30/31, 36/37, 44/45 and 51 also have their Coin name / Link swapped
# or replaced, in addition to new Price/Volume figures.

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.31%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.64%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.309.27"
$ws.Range("E37").Value = "  -5.29%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.949"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.57%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0610"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.70%  "
